$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 30 (shifts old rows 31-47 down to 33-49)
$ws.Rows.Item(31).Insert()
$ws.Rows.Item(31).Insert()

# Update row 30 (Pseudotime -> Pseudotime_1) - only changed columns
$ws.Range('A30').Value = 'Pseudotime_1'
$ws.Range('B30').Value = '  3.76 (2.36)  '
$ws.Range('C30').Value = '  8.25 (1.08)  '
$ws.Range('D30').Value = '  11.2 (0.97)  '
$ws.Range('G30').Value = '  0.40 (0.29)  '
$ws.Range('H30').Value = ' 14.8 (0.99)  '

# Fill new row 31 (Pseudotime_2)
$ws.Range('A31').Value = 'Pseudotime_2'
$ws.Range('B31').Value = '  3.50 (1.53)  '
$ws.Range('C31').Value = '  8.35 (1.12)  '
$ws.Range('D31').Value = '  9.64 (0.22)  '
$ws.Range('E31').Value = '  13.1 (0.48)  '
$ws.Range('F31').Value = '    . (.)     '
$ws.Range('G31').Value = '  0.52 (0.29)  '
$ws.Range('H31').Value = '    . (.)     '
$ws.Range('I31').Value = ' <0.001  '

# Fill new row 32 (Pseudotime_3)
$ws.Range('A32').Value = 'Pseudotime_3'
$ws.Range('B32').Value = '  4.09 (1.87)  '
$ws.Range('C32').Value = '     . (.)     '
$ws.Range('D32').Value = '  8.51 (0.23)  '
$ws.Range('E32').Value = '     . (.)     '
$ws.Range('F32').Value = ' 22.1 (0.04)  '
$ws.Range('G32').Value = '  0.40 (0.30)  '
$ws.Range('H32').Value = ' 9.99 (0.17)  '
$ws.Range('I32').Value = ' <0.001  '
